$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Add header label in J2 and format J3:J5 as dates (mmm-yy, numFmtId 17)
$ws1.Range("J2").Value = "Eclipses fechas"

$ws1.Range("J3:J5").NumberFormat = "mmm-yy"

$ws1.Range("J3").Value = 43770
$ws1.Range("J4").Value = 42795
$ws1.Range("J5").Value = 41821

# Widen columns J and K so the stored width equals 16
$ws1.Columns.Item(10).ColumnWidth = 15.17
$ws1.Columns.Item(11).ColumnWidth = 15.17

# Set selection on Sheet1 to J4:J5 and make Sheet1 the active (tabSelected) sheet
$ws1.Range("J4:J5").Select()
$ws1.Activate()

$wb.Save()
